# Fill in the "Day 6" (G column) sleep-diary answers for the third weekly
# block (rows 44-57) that were left blank, mirroring the F column layout.
# Numeric / time cells use Value2 so the stored number matches Excel's
# serial-date representation exactly; the two "物质使用" answers reuse the
# existing "无" shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wake-up / get-up / bed / lights-off times -> stored as day-fraction serials,
# formatted the same way the neighbouring F column already is (h:mm).
$ws.Range("G44:G47").NumberFormat = "h:mm"

$ws.Range("G44").Value2 = 0.32083333333333336   # 7:42
$ws.Range("G45").Value2 = 0.33333333333333331   # 8:00
$ws.Range("G46").Value2 = 0.96180555555555558   # 23:05
$ws.Range("G47").Value2 = 0.96527777777777779   # 23:10

# Minutes to fall asleep / times woken / minutes awake / minutes slept
$ws.Range("G48").Value2 = 10
$ws.Range("G49").Value2 = 0
$ws.Range("G50").Value2 = 0
$ws.Range("G51").Value2 = 510

# Substance use before bed
$ws.Range("G52").Value2 = "无"

# Electronics use (minutes), physical tension, mental tension, sleep quality
$ws.Range("G53").Value2 = 15
$ws.Range("G54").Value2 = 4
$ws.Range("G55").Value2 = 3
$ws.Range("G56").Value2 = 3

# Daytime nap
$ws.Range("G57").Value2 = "无"

# Leave the selection on the last cell that was filled in, like the author did.
$ws.Range("G57").Select() | Out-Null
